$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Here our modifications 5 for images" -> "...6 for images"
# ------------------------------------------------------------------
$d.Content.Find.Execute("Here our modifications 5 for images", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Here our modifications 6 for images", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Replace the red/Georgia "Image 1:" heading paragraph (which carries
#    the ACE_Image1 bookmark) with three plain Calibri paragraphs:
#       "Image1"
#       "Image 1"
#       "Image 1:"
#    (no bookmark survives).
# ------------------------------------------------------------------
$headingPara = $d.Paragraphs(2)
$headingPara.Range.Delete()

$anchor = $d.Paragraphs(1).Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$d.Paragraphs(2).Range.Text = "Image1"

$anchor = $d.Paragraphs(2).Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$d.Paragraphs(3).Range.Text = "Image 1"

$anchor = $d.Paragraphs(3).Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$d.Paragraphs(4).Range.Text = "Image 1:"

# ------------------------------------------------------------------
# 3) Move "ACE_Image1" + " here" runs up into the previously empty
#    Segoe UI paragraph, then drop the now-redundant paragraph that
#    used to hold that text.
# ------------------------------------------------------------------
$srcRange = $d.Paragraphs(6).Range
$srcRange.Copy()
$destRange = $d.Paragraphs(5).Range
$destRange.Collapse(1)
$destRange.Paste()
$d.Paragraphs(6).Range.Delete()

# ------------------------------------------------------------------
# 4) Drop the stray <w:lastRenderedPageBreak/> before
#    "i hope you are alright" by rewriting the run's text in place
#    (keeps the run's formatting, loses the page-break marker).
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*i hope you are alright*") {
        $p.Range.Text = "i hope you are alright"
    }
}
